# Resize the single-column "representative/contractor/address" signature
# table: its only column shrinks from 6242 dxa (312.1 pt) to 3203 dxa
# (160.15 pt). Setting Column.Width updates both the <w:tblGrid>'s
# <w:gridCol> and every cell's <w:tcW> in that column to match.

$d = $word.ActiveDocument

$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Columns.Count -eq 1 -and $t.Range.Text -like "*{representative}*" -and $t.Range.Text -like "*{contractor}*") {
        $targetTable = $t
        break
    }
}

if ($targetTable -ne $null) {
    $targetTable.Columns.Item(1).Width = 160.15
}
